# Fruta / hortaliza, semanal
# Insert a new weekly price observation (fecha 44753 = 2022-07-11) for
# "Terminal La Palmera de La Serena - Plátano" as the new first data row
# (rows shift down by 3; dimension grows from T729 to T732).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows right above the current row 700 (pushes the
# existing 700-729 block down to 703-732, keeping all their data intact).
$ws.Range("A700:A702").EntireRow.Insert()

# Data for the 3 new rows (Pintón / Primera Maduro / Primera Pintón),
# following the same layout used throughout the sheet for this market.
$fecha   = 44753
$calidad = @("Pintón", "Primera Maduro", "Primera Pintón")
$volumen = @(80, 120, 120)
$precio  = @(21000, 23000, 24000)
$precioKg = @(1050, 1150, 1200)

for ($i = 0; $i -lt 3; $i++) {
    $r = 700 + $i

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108006
    $ws.Cells.Item($r, 10).Value = "Plátano"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $calidad[$i]
    $ws.Cells.Item($r, 13).Value = $volumen[$i]
    $ws.Cells.Item($r, 14).Value = $precio[$i]
    $ws.Cells.Item($r, 15).Value = $precio[$i]
    $ws.Cells.Item($r, 16).Value = $precio[$i]
    $ws.Cells.Item($r, 17).Value = '$/caja 20 kilos'
    $ws.Cells.Item($r, 18).Value = "Ecuador"
    $ws.Cells.Item($r, 19).Value = $precioKg[$i]
    $ws.Cells.Item($r, 20).Value = 20
}
